# Add a new "Postgres LOB (disabled 2nd level hibernate cache, lz4 compression)"
# row to the timings table, right above the existing "Postgres JSONB
# (disabled 2nd level hibernate cache)" row (i.e. as the new row 7, pushing
# every row below it down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at row 7; everything from the old row 7 downward shifts
# down by one (matches the diff: old row7->row8, old row8->row9, ...).
$ws.Rows.Item(7).Insert()

# The new row 7 has no formatting yet - clone it from row 8 (which now holds
# what used to be row 7's data/format) so the new row matches the table's
# look (styles s="4"/"5"/"1").
$ws.Range("A8:H8").Copy()
$ws.Range("A7:H7").PasteSpecial(-4122)

# Fill in the new row's data.
$ws.Range("A7").Value = "Postgres LOB (disabled 2nd level hibernate cache, lz4 compression)"
$ws.Range("B7").Value = 115
$ws.Range("C7").Value = 1257
$ws.Range("D7").Value = 2024
$ws.Range("E7").Value = 152
$ws.Range("F7").Value = 1432
$ws.Range("G7").Value = 2195
$ws.Range("H7").Value = "local docker"

# Match the author's final selection position.
$ws.Range("C22").Select()
